$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.043.41'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '1.800.16'
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '307.69'
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = '0.4212'
$ws.Range('E7').Value = '  -2.18%  '
$ws.Range('D8').Value = '0.3596'
$ws.Range('E8').Value = '  -2.64%  '
$ws.Range('D9').Value = '0.07268'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').Value = '0.8458'
$ws.Range('E10').Value = '  -3.38%  '
$ws.Range('D11').Value = '20.27'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('D12').Value = '1.816.53'
$ws.Range('E12').Value = '  -3.96%  '
$ws.Range('D13').Value = '5.297'
$ws.Range('E13').Value = '  -3.10%  '
$ws.Range('D14').Value = '6.380'
$ws.Range('E14').Value = '  -3.19%  '
$ws.Range('D15').Value = '0.06774'
$ws.Range('E15').Value = '  -2.41%  '
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '80.55'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '0.000008757'
$ws.Range('E18').Value = '  -3.28%  '
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('E20').Value = '  -3.23%  '
$ws.Range('D21').Value = '27.310.41'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = '5.078'
$ws.Range('D23').Value = '11.03'
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').Value = '2.084.55'
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('E25').Value = '  -3.41%  '
$ws.Range('D26').Value = '153.35'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').Value = '18.16'
$ws.Range('E27').Value = '  -4.21%  '
$ws.Range('D28').Value = '5.035'
$ws.Range('E28').Value = '  -5.29%  '
$ws.Range('D29').Value = '113.42'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').Value = '1.658'
$ws.Range('E30').Value = '  -11.36%  '
$ws.Range('D31').Value = '0.08993'
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').Value = '0.7342'
$ws.Range('E32').Value = '  -6.49%  '
$ws.Range('D33').Value = '2.861'
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('D34').Value = '4.348'
$ws.Range('E34').Value = '  -5.50%  '
$ws.Range('E35').Value = '  -5.77%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '1.083'
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('D38').Value = '0.05160'
$ws.Range('E38').Value = '  -5.14%  '
$ws.Range('D39').Value = '0.01907'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.4991'
$ws.Range('E40').Value = '  -3.39%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.1634'
$ws.Range('E41').Value = '  -3.47%  '
$ws.Range('D42').Value = '2.644'
$ws.Range('E42').Value = '  -6.84%  '
$ws.Range('D43').Value = '8.083'
$ws.Range('E43').Value = '  -6.41%  '
$ws.Range('D44').Value = '5.984'
$ws.Range('E44').Value = '  -11.83%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '10.31'
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '105.44'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('E48').Value = '  -3.43%  '
$ws.Range('D49').Value = '0.4538'
$ws.Range('E49').Value = '  -5.31%  '
$ws.Range('D50').Value = '1.605'
$ws.Range('E50').Value = '  -3.36%  '
$ws.Range('D51').Value = '1.739'
$ws.Range('E51').Value = '  -5.61%  '
